$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G11 (quarterly indexing bug-fix: recomputed error-table stats)
$data = @(
    @(0.4414089885834451, 2.069638249730092,  13.91237005874871, 3.729928961622286, 3.748613346451448, 42),
    @(0.6527988889848677, 2.704908097310158,  26.48805711061894, 5.146654943807573, 5.168506407966236, 41),
    @(0.4956225713272083, 3.16172899281031,   28.92945162613211, 5.378610566506197, 5.423955401537688, 40),
    @(0.7736103189250133, 3.244499711943635,  27.52530586389647, 5.246456505480291, 5.25694148770441,  39),
    @(0.4595302257764673, 3.50693379280404,   35.78621063427334, 5.982157690522153, 6.04454549522689,  38),
    @(0.8473814559220318, 3.362369025740158,  31.85202081669558, 5.643759457728118, 5.656747993276202, 37),
    @(0.3348332891986572, 3.296274014622933,  29.66640705715701, 5.446687714304631, 5.513501801737392, 36),
    @(0.5287725799852153, 2.929081919830375,  28.22038967316498, 5.312286670838178, 5.363075344464566, 35),
    @(0.2864215739762718, 3.411304146226604,  30.19520442758122, 5.495016326416257, 5.570070675043795, 34),
    @(0.5683593138408569, 3.343257824933541,  34.92555082416861, 5.909784329750842, 5.973595836522446, 33)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
